$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "'Bitcoin"
$ws.Range("C2").Value = "'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "'68.884.38"
$ws.Range("E2").Value = "'  +1.94%  "
$ws.Range("B2:E2").Style = "Normal"
$ws.Range("B3").Value = "'Ethereum"
$ws.Range("C3").Value = "'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "'3.297.92"
$ws.Range("E3").Value = "'  +1.38%  "
$ws.Range("B3:E3").Style = "Normal"
$ws.Range("B4").Value = "'TetherUSD"
$ws.Range("C4").Value = "'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("B4:E4").Style = "Normal"
$ws.Range("B5").Value = "'BNB"
$ws.Range("C5").Value = "'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'586.70"
$ws.Range("E5").Value = "'  +1.46%  "
$ws.Range("B5:E5").Style = "Normal"
$ws.Range("B6").Value = "'Solana"
$ws.Range("C6").Value = "'https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").Value = "'183.69"
$ws.Range("E6").Value = "'  +1.00%  "
$ws.Range("B6:E6").Style = "Normal"
$ws.Range("B7").Value = "'USDC"
$ws.Range("C7").Value = "'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("B7:E7").Style = "Normal"
$ws.Range("B8").Value = "'XRP"
$ws.Range("C8").Value = "'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.599"
$ws.Range("E8").Value = "'  +1.61%  "
$ws.Range("B8:E8").Style = "Normal"
$ws.Range("B9").Value = "'Dogecoin"
$ws.Range("C9").Value = "'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.135"
$ws.Range("E9").Value = "'  +3.75%  "
$ws.Range("B9:E9").Style = "Normal"
$ws.Range("B10").Value = "'Toncoin"
$ws.Range("C10").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value = "'6.68"
$ws.Range("E10").Value = "'  -1.44%  "
$ws.Range("B10:E10").Style = "Normal"
$ws.Range("B11").Value = "'Cardano"
$ws.Range("C11").Value = "'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").Value = "'0.423"
$ws.Range("E11").Value = "'  +2.24%  "
$ws.Range("B11:E11").Style = "Normal"
$ws.Range("B12").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "'3.874.77"
$ws.Range("E12").Value = "'  +1.67%  "
$ws.Range("B12:E12").Style = "Normal"
$ws.Range("B13").Value = "'TRON"
$ws.Range("C13").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.137"
$ws.Range("E13").Value = "'  -0.23%  "
$ws.Range("B13:E13").Style = "Normal"
$ws.Range("B14").Value = "'Avalanche"
$ws.Range("C14").Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'29.08"
$ws.Range("E14").Value = "'  +3.35%  "
$ws.Range("B14:E14").Style = "Normal"
$ws.Range("B15").Value = "'WrappedBTC"
$ws.Range("C15").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "'68.871.39"
$ws.Range("E15").Value = "'  +1.97%  "
$ws.Range("B15:E15").Style = "Normal"
$ws.Range("B16").Value = "'ShibaInu"
$ws.Range("C16").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000172"
$ws.Range("E16").Value = "'  +2.72%  "
$ws.Range("B16:E16").Style = "Normal"
$ws.Range("B17").Value = "'WrappedEther"
$ws.Range("C17").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'3.326.20"
$ws.Range("E17").Value = "'  +2.33%  "
$ws.Range("B17:E17").Style = "Normal"
$ws.Range("B18").Value = "'Polkadot"
$ws.Range("C18").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'5.85"
$ws.Range("E18").Value = "'  +0.64%  "
$ws.Range("B18:E18").Style = "Normal"
$ws.Range("B19").Value = "'Chainlink"
$ws.Range("C19").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'13.66"
$ws.Range("E19").Value = "'  +1.31%  "
$ws.Range("B19:E19").Style = "Normal"
$ws.Range("B20").Value = "'BitcoinCash"
$ws.Range("C20").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'394.84"
$ws.Range("E20").Value = "'  +5.21%  "
$ws.Range("B20:E20").Style = "Normal"
$ws.Range("B21").Value = "'Uniswap"
$ws.Range("C21").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'7.76"
$ws.Range("E21").Value = "'  +1.99%  "
$ws.Range("B21:E21").Style = "Normal"
$ws.Range("B22").Value = "'Litecoin"
$ws.Range("C22").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").Value = "'71.83"
$ws.Range("E22").Value = "'  +0.77%  "
$ws.Range("B22:E22").Style = "Normal"
$ws.Range("B23").Value = "'Dai"
$ws.Range("C23").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "'  -0.23%  "
$ws.Range("B23:E23").Style = "Normal"
$ws.Range("B24").Value = "'PEPE"
$ws.Range("C24").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D24").Value = "'0.0000122"
$ws.Range("E24").Value = "'  +2.18%  "
$ws.Range("B24:E24").Style = "Normal"
$ws.Range("B25").Value = "'Polygon"
$ws.Range("C25").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").Value = "'0.518"
$ws.Range("E25").Value = "'  +1.38%  "
$ws.Range("B25:E25").Style = "Normal"
$ws.Range("B26").Value = "'Kaspa"
$ws.Range("C26").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "'0.189"
$ws.Range("E26").Value = "'  +4.34%  "
$ws.Range("B26:E26").Style = "Normal"
$ws.Range("B27").Value = "'InternetComputer(DFINITY)"
$ws.Range("C27").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "'9.73"
$ws.Range("E27").Value = "'  +0.68%  "
$ws.Range("B27:E27").Style = "Normal"
$ws.Range("B28").Value = "'Binance-PegBSC-USD"
$ws.Range("C28").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "'  -0.32%  "
$ws.Range("B28:E28").Style = "Normal"
$ws.Range("B29").Value = "'NEARProtocol"
$ws.Range("C29").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").Value = "'5.77"
$ws.Range("E29").Value = "'  +1.12%  "
$ws.Range("B29:E29").Style = "Normal"
$ws.Range("B30").Value = "'PancakeSwap"
$ws.Range("C30").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.99"
$ws.Range("E30").Value = "'  +0.98%  "
$ws.Range("B30:E30").Style = "Normal"
$ws.Range("B31").Value = "'EthereumClassic"
$ws.Range("C31").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'23.18"
$ws.Range("E31").Value = "'  +2.45%  "
$ws.Range("B31:E31").Style = "Normal"
$ws.Range("B32").Value = "'Fetch.AI"
$ws.Range("C32").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.31"
$ws.Range("E32").Value = "'  +2.94%  "
$ws.Range("B32:E32").Style = "Normal"
$ws.Range("B33").Value = "'Aptos"
$ws.Range("C33").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").Value = "'7.19"
$ws.Range("E33").Value = "'  +4.65%  "
$ws.Range("B33:E33").Style = "Normal"
$ws.Range("B34").Value = "'USDe"
$ws.Range("C34").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "'  +0.05%  "
$ws.Range("B34:E34").Style = "Normal"
$ws.Range("B35").Value = "'ImmutableX"
$ws.Range("C35").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'1.54"
$ws.Range("E35").Value = "'  +3.01%  "
$ws.Range("B35:E35").Style = "Normal"
$ws.Range("B36").Value = "'Monero"
$ws.Range("C36").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "'163.51"
$ws.Range("E36").Value = "'  +0.57%  "
$ws.Range("B36:E36").Style = "Normal"
$ws.Range("B37").Value = "'Stacks"
$ws.Range("C37").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").Value = "'1.90"
$ws.Range("E37").Value = "'  +1.86%  "
$ws.Range("B37:E37").Style = "Normal"
$ws.Range("B38").Value = "'Mantle"
$ws.Range("C38").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").Value = "'0.838"
$ws.Range("E38").Value = "'  -2.02%  "
$ws.Range("B38:E38").Style = "Normal"
$ws.Range("B39").Value = "'Filecoin"
$ws.Range("C39").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").Value = "'4.62"
$ws.Range("E39").Value = "'  +3.74%  "
$ws.Range("B39:E39").Style = "Normal"
$ws.Range("B40").Value = "'EnergySwap"
$ws.Range("C40").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "'26.44"
$ws.Range("E40").Value = "'  -1.16%  "
$ws.Range("B40:E40").Style = "Normal"
$ws.Range("B41").Value = "'RenderToken"
$ws.Range("C41").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'6.57"
$ws.Range("E41").Value = "'  -3.44%  "
$ws.Range("B41:E41").Style = "Normal"
$ws.Range("B42").Value = "'dogwifhat"
$ws.Range("C42").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.59"
$ws.Range("E42").Value = "'  -0.44%  "
$ws.Range("B42:E42").Style = "Normal"
$ws.Range("B43").Value = "'OKB"
$ws.Range("C43").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "'41.58"
$ws.Range("E43").Value = "'  +2.49%  "
$ws.Range("B43:E43").Style = "Normal"
$ws.Range("B44").Value = "'Hedera"
$ws.Range("C44").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").Value = "'0.0693"
$ws.Range("E44").Value = "'  +2.72%  "
$ws.Range("B44:E44").Style = "Normal"
$ws.Range("B45").Value = "'Bittensor"
$ws.Range("C45").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "'344.53"
$ws.Range("E45").Value = "'  -5.02%  "
$ws.Range("B45:E45").Style = "Normal"
$ws.Range("B46").Value = "'InjectiveProtocol"
$ws.Range("C46").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'25.23"
$ws.Range("E46").Value = "'  -1.36%  "
$ws.Range("B46:E46").Style = "Normal"
$ws.Range("B47").Value = "'Maker"
$ws.Range("C47").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "'2.621.51"
$ws.Range("E47").Value = "'  -4.31%  "
$ws.Range("B47:E47").Style = "Normal"
$ws.Range("B48").Value = "'VeChain"
$ws.Range("C48").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0284"
$ws.Range("E48").Value = "'  +2.00%  "
$ws.Range("B48:E48").Style = "Normal"
$ws.Range("B49").Value = "'Arweave"
$ws.Range("C49").Value = "'https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").Value = "'32.25"
$ws.Range("E49").Value = "'  +4.37%  "
$ws.Range("B49:E49").Style = "Normal"
$ws.Range("B50").Value = "'Cosmos"
$ws.Range("C50").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'6.33"
$ws.Range("E50").Value = "'  +3.46%  "
$ws.Range("B50:E50").Style = "Normal"
$ws.Range("B51").Value = "'Stellar"
$ws.Range("C51").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.102"
$ws.Range("E51").Value = "'  +0.35%  "
$ws.Range("B51:E51").Style = "Normal"
